$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-like string as plain text (shared string) without
# letting Excel auto-convert it into a date serial number + date format.
# We stage the text in a scratch cell via a formula (so it is typed as
# text), copy it, and paste-special (values only) into the target cell.
# This avoids introducing any cell style, matching how the original rows
# store their DATE column as plain text.
function Set-TextValue($cellRef, $text) {
    $ws.Range("Z1").Formula = '="' + $text + '"'
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("Z1").ClearContents()
    $excel.CutCopyMode = $false
}

# Row 48: 2022-07-06 / WEDNESDAY
Set-TextValue "A48" "2022-07-06"
$ws.Range("B48").Value = "WEDNESDAY"
$ws.Range("C48").Value = 6.100000000000001
$ws.Range("D48").Value = 152.50000000000003
$ws.Range("E48").Value = 610.0000000000001

# Row 49: 2022-07-07 / THURSDAY
Set-TextValue "A49" "2022-07-07"
$ws.Range("B49").Value = "THURSDAY"
$ws.Range("C49").Value = 36.45
$ws.Range("D49").Value = 911.2500000000001
$ws.Range("E49").Value = 3645.0000000000005
